# "good events, not on screen yet"
#
# 1. Fix a typo in the BAD sheet (B6): "pug up" -> "plug up".
# 2. Row 10 on BAD currently holds answers that actually belong to a
#    different (mis-pasted) event; row 11 holds the correct answers for
#    its own event label but they were duplicated into row 10 as well.
#    Move row 11's answers up into row 10 and blank out row 11's answers
#    (keeping its own event label in column A).
# 3. Switch the active sheet/tab from BAD to GOOD, and update each
#    sheet's remembered selection/scroll position.

$wb = $excel.ActiveWorkbook
$bad = $wb.Worksheets.Item("BAD")
$good = $wb.Worksheets.Item("GOOD")

# --- 1. Typo fix -----------------------------------------------------
$bad.Range("B6").Value2 = 'You simply plug up the volcano, and the heat that is coming from the mountain fuel your followers famous "sfefse" springs'

# --- 2. Move row 11's outcomes up into row 10, blank row 11 ----------
$bad.Range("B10:I10").Value2 = $bad.Range("B11:I11").Value2
$bad.Range("B11:I11").ClearContents()

$bad.Rows.Item(10).RowHeight = 105
$bad.Rows.Item(11).AutoFit()

# --- 3. Selection / active tab ----------------------------------------
$bad.Range("B10").Select() | Out-Null

$good.Activate() | Out-Null
$good.Range("B2").Select() | Out-Null
